$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.01135533333333333
$ws.Range("H2").Value = 0.034066
$ws.Range("M2").Value = 61.04160633333334
$ws.Range("N2").Value = 183.124819
$ws.Range("O2").Value = 0.2043613460574534
$ws.Range("P2").Value = 0.2043613460574534
$ws.Range("Q2").Value = 0.6931477871171112
$ws.Range("R2").Value = 6.238330084054001
$ws.Range("S2").Value = 0.2043613460574534
$ws.Range("T2").Value = 0.2043613460574534

# Row 3
$ws.Range("G3").Value = 0.01135533333333333
$ws.Range("H3").Value = 0.034066
$ws.Range("O3").Value = 0.3559304658284363
$ws.Range("P3").Value = 0.3559304658284363
$ws.Range("Q3").Value = 1.207236199585333
$ws.Range("R3").Value = 10.865125796268
$ws.Range("S3").Value = 0.3559304658284363
$ws.Range("T3").Value = 0.3559304658284363

# Row 4
$ws.Range("G4").Value = 0.01135533333333333
$ws.Range("H4").Value = 0.034066
$ws.Range("M4").Value = 131.3384093333333
$ws.Range("N4").Value = 394.015228
$ws.Range("O4").Value = 0.4397081881141102
$ws.Range("P4").Value = 0.4397081881141103
$ws.Range("Q4").Value = 1.491391417449778
$ws.Range("R4").Value = 13.422522757048
$ws.Range("S4").Value = 0.4397081881141102
$ws.Range("T4").Value = 0.4397081881141103
